$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = -10.85
$ws.Range("C18").Value = -11.62
$ws.Range("C20").Value = -12.298
$ws.Range("C27").Value = -12.951
$ws.Range("C35").Value = -12.197
$ws.Range("C69").Value = -10.537
$ws.Range("C76").Value = -13.201
$ws.Range("C78").Value = -12.5
$ws.Range("C82").Value = -11.991
$ws.Range("C83").Value = -13.262
$ws.Range("C93").Value = -11.642
